$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'silicone knee pads'
$ws.Range('A2').Value = 'skins compression men'
$ws.Range('A3').Value = 'sliding knee pads baseball'
$ws.Range('A4').Value = 'small football knee pads'
$ws.Range('A5').Value = 'snowboarding knee pad'
$ws.Range('A6').Value = 'soccer knee pad'
$ws.Range('A7').Value = 'soccer knee protector'
$ws.Range('A8').Value = 'soccer pants for boys'
$ws.Range('A9').Value = 'spandex pants for men'
$ws.Range('A10').Value = 'spats bjj'
$ws.Range('A11').Value = 'sports basketball pants'
$ws.Range('A12').Value = 'sports leggings men'
$ws.Range('A13').Value = 'squat protector'
$ws.Range('A14').Value = 'strip pants men'
$ws.Range('A15').Value = 'tactical pants with knee pads'
$ws.Range('A16').Value = 'thermal baselayer men'
$ws.Range('A17').Value = 'thermal compression pants men'
$ws.Range('A18').Value = 'thermal leggings men'
$ws.Range('A19').Value = 'thermal winter pants'
$ws.Range('A20').Value = 'thin knee pads'
$ws.Range('A21').Value = 'tight pant'
$ws.Range('A22').Value = 'tights black'
$ws.Range('A23').Value = 'tights for men'
$ws.Range('A24').Value = 'under amour fleece leggings'
$ws.Range('A25').Value = 'under armour leggings'
$ws.Range('A26').Value = 'underarmor thermal pants mens'
$ws.Range('A27').Value = 'volleyball clothes men'
$ws.Range('A28').Value = 'volleyball compression knee pads'
$ws.Range('A29').Value = 'volleyball knee pads adult'
$ws.Range('A30').Value = 'volleyball knee pads youth'
$ws.Range('A31').Value = 'warm compression pants mens'
$ws.Range('A32').Value = 'weightlifting equipment'
$ws.Range('A33').Value = 'white basketball knee pads'
$ws.Range('A34').Value = 'white tights mens basketball'
$ws.Range('A35').Value = 'womens basketball knee pads'
$ws.Range('A36').Value = 'workout pads for hands'
$ws.Range('A37').Value = 'workout squat pad'
$ws.Range('A38').Value = 'wrestling knee pads adult'
$ws.Range('A39').Value = 'wrestling tights boys'
$ws.Range('A40').Value = 'xl volleyball knee pads'
$ws.Range('A41').Value = 'youth basketball knee'
$ws.Range('A42').Value = 'youth basketball leggings for boys'
$ws.Range('A43').Value = 'youth basketball pants girls'
$ws.Range('A44').Value = 'youth football knee pads'
$ws.Range('A45').Value = 'youth knee pads for skating'
$ws.Range('A46').Value = 'youth knee pads mountain bike'
$ws.Range('A47').Value = 'youth soccer gear'
$ws.Range('A48').Value = 'youth sports tights'
$ws.Range('A49').Value = 'mens running tights capri'
$ws.Range('A50').Value = 'mens basketball pants tall'
$ws.Range('A51').Value = 'volleyball hip protectors'
$ws.Range('A52').Value = 'knee compression basketball'
$ws.Range('A53').Value = 'sports compression pants'
$ws.Range('A54').Value = 'basketball leggings youth'
$ws.Range('A55').Value = 'knee pads for running'
$ws.Range('A56').Value = 'wrestling pants for men'
$ws.Range('A57').Value = 'male compression pants'
$ws.Range('A58').Value = 'cheap compression pants men'
$ws.Range('A59').Value = 'men compression pants pack'
$ws.Range('A60').Value = 'knee pad men'
$ws.Range('A61').Value = 'volleyball knee pads for youth'
$ws.Range('A62').Value = 'leggings for man'
$ws.Range('A63').Value = 'volleyball knee pads xl'
$ws.Range('A64').Value = 'spandex leggings men'
$ws.Range('A65').Value = 'knees compression'
$ws.Range('A66').Value = 'athletic boys pants'
$ws.Range('A67').Value = 'soccer knee pads'
$ws.Range('A68').Value = 'adults knee pads'
$ws.Range('A69').Value = 'leg compression leggings'
$ws.Range('A70').Value = 'big and tall tights for men'
$ws.Range('A71').Value = 'wrestling mens apparel'
$ws.Range('A72').Value = 'mens hiking leggings'
$ws.Range('A73').Value = 'softball pants mens'
$ws.Range('A74').Value = 'running compression knee'
$ws.Range('A75').Value = 'compression pads for basketball'
$ws.Range('A76').Value = 'gym tights'
$ws.Range('A77').Value = 'boys soccer leggings'
$ws.Range('A78').Value = 'hockey pants men'
$ws.Range('A79').Value = 'knee pads for volleyball for men'
$ws.Range('A80').Value = 'baseball pants men'
$ws.Range('A81').Value = 'mens running pants'
$ws.Range('A82').Value = 'sports pants for men'
$ws.Range('A83').Value = 'boys snowboarding pants'
$ws.Range('A84').Value = 'clothing protector'
$ws.Range('A85').Value = 'baseball pants for youth'
$ws.Range('A86').Value = 'mens gym pants'
$ws.Range('A87').Value = 'football pants adult'
$ws.Range('A88').Value = 'knee pads sports'
$ws.Range('A89').Value = 'protective knee pads'
$ws.Range('A90').Value = 'girls volleyball pads'
$ws.Range('A91').Value = 'pants youth'
$ws.Range('A92').Value = 'big boys compression leggings'
$ws.Range('A93').Value = 'basketball spandex'
$ws.Range('A94').Value = 'mens basketball clothing'
$ws.Range('A95').Value = 'boys tights for sports'
$ws.Range('A96').Value = 'leggings for cycling'
$ws.Range('A97').Value = 'hiking capri men'
$ws.Range('A98').Value = 'baseball pants men black'
$ws.Range('A99').Value = 'football tights youth'
$ws.Range('A100').Value = 'mens 3/4 tights'
